$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.03397833333333333
$ws.Range("H2").Value = 0.101935
$ws.Range("I2").Value = 0.03987224921182536
$ws.Range("J2").Value = 0.03987224921182535
$ws.Range("M2").Value = 15.24491733333333
$ws.Range("N2").Value = 45.73475199999999
$ws.Range("O2").Value = 0.4831257321597052
$ws.Range("P2").Value = 0.4831257321597052
$ws.Range("Q2").Value = 0.517996882791111
$ws.Range("R2").Value = 4.661971945119999
$ws.Range("S2").Value = 0.01926330959331736
$ws.Range("T2").Value = 0.01926330959331735

$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.03397833333333333
$ws.Range("H3").Value = 0.101935
$ws.Range("I3").Value = 0.03987224921182536
$ws.Range("J3").Value = 0.03987224921182535
$ws.Range("O3").Value = 0.327710667227878
$ws.Range("P3").Value = 0.327710667227878
$ws.Range("Q3").Value = 0.35136423664
$ws.Range("R3").Value = 3.16227812976
$ws.Range("S3").Value = 0.01306656139308352
$ws.Range("T3").Value = 0.01306656139308352

$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.03397833333333333
$ws.Range("H4").Value = 0.101935
$ws.Range("I4").Value = 0.03987224921182536
$ws.Range("J4").Value = 0.03987224921182535
$ws.Range("M4").Value = 5.969012333333333
$ws.Range("N4").Value = 17.907037
$ws.Range("O4").Value = 0.1891636006124168
$ws.Range("P4").Value = 0.1891636006124168
$ws.Range("Q4").Value = 0.2028170907327778
$ws.Range("R4").Value = 1.825353816595
$ws.Range("S4").Value = 0.007542378225424484
$ws.Range("T4").Value = 0.007542378225424483

$ws.Range("G5").Value = 0.7475459999999999
$ws.Range("I5").Value = 0.8772160811096247
$ws.Range("J5").Value = 0.8772160811096247
$ws.Range("M5").Value = 15.24491733333333
$ws.Range("N5").Value = 45.73475199999999
$ws.Range("O5").Value = 0.4831257321597052
$ws.Range("P5").Value = 0.4831257321597052
$ws.Range("Q5").Value = 11.396276972864
$ws.Range("R5").Value = 102.566492755776
$ws.Range("S5").Value = 0.4238056614483548
$ws.Range("T5").Value = 0.4238056614483547

$ws.Range("G6").Value = 0.7475459999999999
$ws.Range("I6").Value = 0.8772160811096247
$ws.Range("J6").Value = 0.8772160811096247
$ws.Range("O6").Value = 0.327710667227878
$ws.Range("P6").Value = 0.327710667227878
$ws.Range("Q6").Value = 7.730247598272
$ws.Range("R6").Value = 69.572228384448
$ws.Range("S6").Value = 0.2874730672434594
$ws.Range("T6").Value = 0.2874730672434594

$ws.Range("G7").Value = 0.7475459999999999
$ws.Range("I7").Value = 0.8772160811096247
$ws.Range("J7").Value = 0.8772160811096247
$ws.Range("M7").Value = 5.969012333333333
$ws.Range("N7").Value = 17.907037
$ws.Range("O7").Value = 0.1891636006124168
$ws.Range("P7").Value = 0.1891636006124168
$ws.Range("Q7").Value = 4.462111293733999
$ws.Range("R7").Value = 40.159001643606
$ws.Range("S7").Value = 0.1659373524178105
$ws.Range("T7").Value = 0.1659373524178105

$ws.Range("G8").Value = 0.07065566666666666
$ws.Range("H8").Value = 0.211967
$ws.Range("I8").Value = 0.08291166967854992
$ws.Range("J8").Value = 0.0829116696785499
$ws.Range("M8").Value = 15.24491733333333
$ws.Range("N8").Value = 45.73475199999999
$ws.Range("O8").Value = 0.4831257321597052
$ws.Range("P8").Value = 0.4831257321597052
$ws.Range("Q8").Value = 1.077139797464889
$ws.Range("R8").Value = 9.694258177183999
$ws.Range("S8").Value = 0.04005676111803306
$ws.Range("T8").Value = 0.04005676111803305

$ws.Range("G9").Value = 0.07065566666666666
$ws.Range("H9").Value = 0.211967
$ws.Range("I9").Value = 0.08291166967854992
$ws.Range("J9").Value = 0.0829116696785499
$ws.Range("O9").Value = 0.327710667227878
$ws.Range("P9").Value = 0.327710667227878
$ws.Range("Q9").Value = 0.730638378848
$ws.Range("R9").Value = 6.575745409632
$ws.Range("S9").Value = 0.02717103859133501
$ws.Range("T9").Value = 0.02717103859133501

$ws.Range("G10").Value = 0.07065566666666666
$ws.Range("H10").Value = 0.211967
$ws.Range("I10").Value = 0.08291166967854992
$ws.Range("J10").Value = 0.0829116696785499
$ws.Range("M10").Value = 5.969012333333333
$ws.Range("N10").Value = 17.907037
$ws.Range("O10").Value = 0.1891636006124168
$ws.Range("P10").Value = 0.1891636006124168
$ws.Range("Q10").Value = 0.4217445457532221
$ws.Range("R10").Value = 3.795700911778999
$ws.Range("S10").Value = 0.01568386996918185
$ws.Range("T10").Value = 0.01568386996918185
